$d = $word.ActiveDocument

# --- Change 1: "subtillem<exp>ent</exp> destrempée" -> "subtilem<exp>ent</exp> destrempée" ---
# This text (with the " destrempée" suffix) is unique in the document, so a scoped
# Find/Replace touches only this occurrence and leaves the other "subtillem...pass" text alone.
$null = $d.Content.Find.Execute("subtillem<exp>ent</exp> destrempée", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "subtilem<exp>ent</exp> destrempée", 2)

# --- Change 2: ". Garde aussi que le gect ne soit point trop large ne" ---
#   -> ". Garde auss" + "y" (distinct run, no explicit color) + " que le gect ne soit point trop large ne"
$findRange = $d.Content
$null = $findRange.Find.Execute("aussi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Force a run boundary right after the "i" we are about to replace, by toggling Bold on the
# trailing text (this causes Word to split the run); we un-bold it again afterwards.
$tailStr = " que le gect ne soit point trop large ne"
$afterRange = $d.Range($findRange.End, $findRange.End + $tailStr.Length)
$afterRange.Bold = $true

# Select just the "i" in "aussi" and replace it with a "y" run that carries no explicit color
# (InsertXML lets us author the run's rPr precisely, producing <w:rPr><w:rtl w:val="0"/></w:rPr>).
$iRange = $d.Range($findRange.End - 1, $findRange.End)
$frag = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t>y</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$iRange.InsertXML($frag)

# Remove the temporary Bold we applied to force the split.
$tailFind = $d.Content
$null = $tailFind.Find.Execute($tailStr, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailFind.Bold = $false
